# Apply the change described by the diff:
#  - Metadata sheet: set "Date" value (B8) to the new timestamp
#  - Metadata sheet: set "Jurisdiction" value (B11) to "FRANCE"
#
# (The shared-string table reshuffling and the corresponding index shifts
#  in the "Include #0" sheet's rows 3/4 are a natural consequence of Excel
#  re-serializing the shared strings table after the cell value change and
#  require no direct edits themselves.)

$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")

$wsMeta.Range("B8").Value = "2025-07-11T12:29:53+00:00"
$wsMeta.Range("B11").Value = "FRANCE"
